# Weekly update: insert a new price record as row 29, pushing the
# existing rows 29..94 down to 30..95 (dimension grows from R94 to R95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 29 (shifts rows 29-94 -> 30-95).
$ws.Rows(29).Insert()

# Populate the new row 29 with the new data point.
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = 100112031
$ws.Cells.Item(29, 7).Value = "Poroto verde"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 35000
$ws.Cells.Item(29, 12).Value = 35000
$ws.Cells.Item(29, 13).Value = 35000
$ws.Cells.Item(29, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(29, 16).Value = 1400
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
